# "study 2 analysis update" - update the H7.x labels on the Study 2
# analysis slide (sldId 264 / 3rd slide) to H5.x.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# TextBox 21 (shape id 22): "H7.1" -> "H5.1"
$s.Shapes.Item(14).TextFrame.TextRange.Text = "H5.1"

# TextBox 23 (shape id 24): "H7.2" -> "H5.2"
# Use Find&Replace on the digit so only the changed character gets a new
# run, matching how the label was actually edited.
$s.Shapes.Item(15).TextFrame.TextRange.Replace("7", "5", 0, 0, 0)

# TextBox 24 (shape id 25): "H7.3" -> "H5.3"
$s.Shapes.Item(16).TextFrame.TextRange.Text = "H5.3"
